$d = $word.ActiveDocument

$d.Content.Find.Execute("2024-05-03 Friday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-05-04 Saturday", 2) | Out-Null
$d.Content.Find.Execute("106×9=", $true, $false, $false, $false, $false, $true, 1, $false, "614×8=", 2) | Out-Null
$d.Content.Find.Execute("848×9=", $true, $false, $false, $false, $false, $true, 1, $false, "359×8=", 2) | Out-Null
$d.Content.Find.Execute("901×7=", $true, $false, $false, $false, $false, $true, 1, $false, "848×3=", 2) | Out-Null
$d.Content.Find.Execute("955×8=", $true, $false, $false, $false, $false, $true, 1, $false, "750×7=", 2) | Out-Null
$d.Content.Find.Execute("515×9=", $true, $false, $false, $false, $false, $true, 1, $false, "766×9=", 2) | Out-Null
$d.Content.Find.Execute("734×7=", $true, $false, $false, $false, $false, $true, 1, $false, "344×6=", 2) | Out-Null
$d.Content.Find.Execute("872×4=", $true, $false, $false, $false, $false, $true, 1, $false, "991×2=", 2) | Out-Null
$d.Content.Find.Execute("428×3=", $true, $false, $false, $false, $false, $true, 1, $false, "837×2=", 2) | Out-Null
$d.Content.Find.Execute("505×6=", $true, $false, $false, $false, $false, $true, 1, $false, "294×5=", 2) | Out-Null
$d.Content.Find.Execute("428×2=", $true, $false, $false, $false, $false, $true, 1, $false, "904×5=", 2) | Out-Null
$d.Content.Find.Execute("423×8=", $true, $false, $false, $false, $false, $true, 1, $false, "255×9=", 2) | Out-Null
$d.Content.Find.Execute("985×8=", $true, $false, $false, $false, $false, $true, 1, $false, "151×6=", 2) | Out-Null
$d.Content.Find.Execute("480×3=", $true, $false, $false, $false, $false, $true, 1, $false, "921×8=", 2) | Out-Null
$d.Content.Find.Execute("286×2=", $true, $false, $false, $false, $false, $true, 1, $false, "953×2=", 2) | Out-Null
$d.Content.Find.Execute("876×3=", $true, $false, $false, $false, $false, $true, 1, $false, "778×2=", 2) | Out-Null
$d.Content.Find.Execute("571×4=", $true, $false, $false, $false, $false, $true, 1, $false, "137×9=", 2) | Out-Null
$d.Content.Find.Execute("534×4=", $true, $false, $false, $false, $false, $true, 1, $false, "473×3=", 2) | Out-Null
$d.Content.Find.Execute("539×7=", $true, $false, $false, $false, $false, $true, 1, $false, "166×3=", 2) | Out-Null
$d.Content.Find.Execute("993×2=", $true, $false, $false, $false, $false, $true, 1, $false, "849×4=", 2) | Out-Null
$d.Content.Find.Execute("465×4=", $true, $false, $false, $false, $false, $true, 1, $false, "506×2=", 2) | Out-Null
$d.Content.Find.Execute("418×7=", $true, $false, $false, $false, $false, $true, 1, $false, "266×5=", 2) | Out-Null
$d.Content.Find.Execute("987×6=", $true, $false, $false, $false, $false, $true, 1, $false, "124×8=", 2) | Out-Null
$d.Content.Find.Execute("296×8=", $true, $false, $false, $false, $false, $true, 1, $false, "166×4=", 2) | Out-Null
$d.Content.Find.Execute("103×3=", $true, $false, $false, $false, $false, $true, 1, $false, "907×9=", 2) | Out-Null
$d.Content.Find.Execute("909×6=", $true, $false, $false, $false, $false, $true, 1, $false, "309×8=", 2) | Out-Null
